$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "CROX STOCK NEWS: CROX Shareholders with Large Losses Should Contact Robbins LLP for Information About the Class Action Lawsuit Against Crocs, Inc.",
    "NTLA Investors Have Opportunity to Lead Intellia Therapeutics, Inc. Securities Fraud Lawsuit",
    "The prosecutor’s resignation letter",
    "ATTOM: Foreclosure Starts Increased 8 Percent in January",
    "How to Request&mdash;and Not Request&mdash;Company Records",
    "Truck stop developer suing city of Monteagle over building permits"
)

$startRow = 134
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $titles[$i]
}
